$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A455").Value = "Buying Opportunity"
$ws.Range("B455").Value = "support Zone"
$ws.Range("C455").Value = "long buildup"
$ws.Range("D455").Value = "Short buildup"
$ws.Range("E455").Value = "FII ENTERING"

$ws.Range("A456").Value = "AHLUCONT"
$ws.Range("B456").Value = "ACL"
$ws.Range("E456").Value = "BATAINDIA"
$ws.Range("F456").Value = 1361.1
$ws.Range("G456").Value = 86.35
$ws.Range("J456").Value = 1373.95

$ws.Range("A457").Value = "NIRAJ"
$ws.Range("B457").Value = "ANDHRAPAP"
$ws.Range("F457").Value = 48.55
$ws.Range("G457").Value = 504.05

$ws.Range("B458").Value = "APOLLOTYRE"
$ws.Range("G458").Value = 471.05

$ws.Range("B459").Value = "ATGL"
$ws.Range("G459").Value = 950.75

$ws.Range("B460").Value = "BALMLAWRIE"
$ws.Range("G460").Value = 276.0

$ws.Range("B461").Value = "BANCOINDIA"
$ws.Range("G461").Value = 615.45

$ws.Range("B462").Value = "BEPL"
$ws.Range("G462").Value = 93.8

$ws.Range("B463").Value = "BPL"
$ws.Range("G463").Value = 98.0

$ws.Range("B464").Value = "CARBORUNIV"
$ws.Range("G464").Value = 1595.0

$ws.Range("B465").Value = "CONSUMBEES"
$ws.Range("G465").Value = 116.93

$ws.Range("B466").Value = "CYBERTECH"
$ws.Range("G466").Value = 145.65

$ws.Range("B467").Value = "DELTAMAGNT"
$ws.Range("G467").Value = 96.75

$ws.Range("B468").Value = "DIAMONDYD"
$ws.Range("G468").Value = 833.6

$ws.Range("B469").Value = "DIVGIITTS"
$ws.Range("G469").Value = 739.75

$ws.Range("B470").Value = "DPWIRES"
$ws.Range("G470").Value = 471.55

$ws.Range("B471").Value = "EKC"
$ws.Range("G471").Value = 126.4

$ws.Range("B472").Value = "EMUDHRA"
$ws.Range("G472").Value = 746.95

$ws.Range("B473").Value = "EPIGRAL"
$ws.Range("G473").Value = 1220.25

$ws.Range("B474").Value = "FACT"
$ws.Range("G474").Value = 693.4

$ws.Range("B475").Value = "GABRIEL"
$ws.Range("G475").Value = 370.1

$ws.Range("B476").Value = "GANESHHOUC"
$ws.Range("G476").Value = 763.65

$ws.Range("B477").Value = "GENCON"
$ws.Range("G477").Value = 39.1

$ws.Range("B478").Value = "GMRP&UI"
$ws.Range("G478").Value = 64.7

$ws.Range("B479").Value = "GOLDIAM"
$ws.Range("G479").Value = 162.85

$ws.Range("B480").Value = "GSFC"
$ws.Range("G480").Value = 218.0

$ws.Range("B481").Value = "GUJALKALI"
$ws.Range("G481").Value = 777.15

$ws.Range("B482").Value = "GULFOILLUB"
$ws.Range("G482").Value = 960.1

$ws.Range("B483").Value = "GULPOLY"
$ws.Range("G483").Value = 181.05

$ws.Range("B484").Value = "HINDMOTORS"
$ws.Range("G484").Value = 34.1

$ws.Range("B485").Value = "INDOAMIN"
$ws.Range("G485").Value = 119.2

$ws.Range("B486").Value = "INFIBEAM"
$ws.Range("G486").Value = 30.6

$ws.Range("B487").Value = "ITDC"
$ws.Range("G487").Value = 651.95

$ws.Range("B488").Value = "JAYSREETEA"
$ws.Range("G488").Value = 98.4

$ws.Range("B489").Value = "KOTARISUG"
$ws.Range("G489").Value = 54.7

$ws.Range("B490").Value = "LUMAXIND"
$ws.Range("G490").Value = 2529.0

$ws.Range("B491").Value = "LXCHEM"
$ws.Range("G491").Value = 243.3

$ws.Range("B492").Value = "MAGNUM"
$ws.Range("G492").Value = 50.75

$ws.Range("B493").Value = "MAXESTATES"
$ws.Range("G493").Value = 366.45

$ws.Range("B494").Value = "MGEL"
$ws.Range("G494").Value = 20.8

$ws.Range("B495").Value = "MICEL"
$ws.Range("G495").Value = 49.05

$ws.Range("B496").Value = "MOLDTECH"
$ws.Range("G496").Value = 238.4

$ws.Range("B497").Value = "MONARCH"
$ws.Range("G497").Value = 528.2

$ws.Range("B498").Value = "MRF"
$ws.Range("G498").Value = 128147.85

$ws.Range("B499").Value = "NAVINIFTY"
$ws.Range("G499").Value = 229.48

$ws.Range("B500").Value = "NAVNETEDUL"
$ws.Range("G500").Value = 150.35

$ws.Range("B501").Value = "NIITLTD"
$ws.Range("G501").Value = 102.9

$ws.Range("B502").Value = "NLCINDIA"
$ws.Range("G502").Value = 221.4

$ws.Range("B503").Value = "OCCL"
$ws.Range("G503").Value = 653.05

$ws.Range("B504").Value = "PATINTLOG"
$ws.Range("G504").Value = 20.95

$ws.Range("B505").Value = "PRICOLLTD"
$ws.Range("G505").Value = 440.1

$ws.Range("B506").Value = "PSPPROJECT"
$ws.Range("G506").Value = 639.5

$ws.Range("B507").Value = "RADIANTCMS"
$ws.Range("G507").Value = 78.75

$ws.Range("B508").Value = "REPRO"
$ws.Range("G508").Value = 732.75

$ws.Range("B509").Value = "RIIL"
$ws.Range("G509").Value = 1214.2

$ws.Range("B510").Value = "RKEC"
$ws.Range("G510").Value = 87.4

$ws.Range("A511").Value = "29/05/2024"

Write-Host "Rows 455-511 written"